$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Switch to Sheet2 (becomes the active/selected tab)
$ws2.Activate()

# Add new header and value for the screenshot interval setting
$ws2.Range("B1").Value = "Screenshot interval (in seconds)"
$ws2.Range("B2").Value = 30

# Widen column B to fit the new header text
$ws2.Columns.Item(2).ColumnWidth = 32

# Leave the selection on B3, as a user would after entering B2's value
$ws2.Range("B3").Select()
